$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1444.5714
$ws.Range("I2").Value = 1421.25
$ws.Range("J2").Value = 1519.2
$ws.Range("K2").Value = 1421.25
$ws.Range("L2").Value = 1519.2
$ws.Range("M2").Value = -1308.25
$ws.Range("N2").Value = -1745.2
$ws.Range("H19").Value = 4321.75
$ws.Range("I19").Value = 8090.125
$ws.Range("K19").Value = 8090.125
$ws.Range("M19").Value = -7915.125
$ws.Range("H62").Value = 130869.625
$ws.Range("I62").Value = 336862.66
$ws.Range("J62").Value = 7273.8
$ws.Range("K62").Value = 336862.66
$ws.Range("L62").Value = 7273.8
$ws.Range("M62").Value = -336238.66
$ws.Range("N62").Value = -8521.799999999999
$ws.Range("H65").Value = 130869.625
$ws.Range("I65").Value = 336862.66
$ws.Range("J65").Value = 7273.8
$ws.Range("K65").Value = 1684313.3
$ws.Range("L65").Value = 36369
$ws.Range("M65").Value = -1681193.3
$ws.Range("N65").Value = -42609
$ws.Range("H96").Value = 182110.73
$ws.Range("J96").Value = 525
$ws.Range("L96").Value = 1575
$ws.Range("N96").Value = -4321
$ws.Range("H112").Value = 1980.75
$ws.Range("J112").Value = 2721
$ws.Range("L112").Value = 8163
$ws.Range("N112").Value = -10379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2722.3
$ws.Range("I2").Value = 2153.5
$ws.Range("K2").Value = 2153.5
$ws.Range("M2").Value = -2040.5
$ws.Range("H32").Value = 160548.8
$ws.Range("I32").Value = 177829.53
$ws.Range("J32").Value = 101794.3
$ws.Range("K32").Value = 177829.53
$ws.Range("L32").Value = 101794.3
$ws.Range("M32").Value = -177542.53
$ws.Range("N32").Value = -102368.3
$ws.Range("H116").Value = 2722.3
$ws.Range("I116").Value = 2153.5
$ws.Range("K116").Value = 2153.5
$ws.Range("M116").Value = 140.5
$ws.Range("H122").Value = 12460.315
$ws.Range("I122").Value = 14473.75
$ws.Range("K122").Value = 43421.25
$ws.Range("M122").Value = -40971.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2722.3
$ws.Range("I3").Value = 2153.5
$ws.Range("K3").Value = 2153.5
$ws.Range("M3").Value = -2039.5
$ws.Range("H82").Value = 19287.691
$ws.Range("J82").Value = 24996.25
$ws.Range("L82").Value = 24996.25
$ws.Range("N82").Value = -25762.25
$ws.Range("H85").Value = 19287.691
$ws.Range("J85").Value = 24996.25
$ws.Range("L85").Value = 24996.25
$ws.Range("N85").Value = -27648.25
$ws.Range("H86").Value = 1970.5
$ws.Range("I86").Value = 1699.8572
$ws.Range("J86").Value = 2602
$ws.Range("K86").Value = 1699.8572
$ws.Range("L86").Value = 2602
$ws.Range("M86").Value = -576.8571999999999
$ws.Range("N86").Value = -4848
$ws.Range("H89").Value = 1970.5
$ws.Range("I89").Value = 1699.8572
$ws.Range("J89").Value = 2602
$ws.Range("K89").Value = 8499.286
$ws.Range("L89").Value = 13010
$ws.Range("M89").Value = -2883.286
$ws.Range("N89").Value = -24242
$ws.Range("H105").Value = 8339667.5
$ws.Range("I105").Value = 11118778
$ws.Range("J105").Value = 2337
$ws.Range("K105").Value = 11118778
$ws.Range("L105").Value = 2337
$ws.Range("M105").Value = -11117031
$ws.Range("N105").Value = -5831
$ws.Range("H134").Value = 2412.8667
$ws.Range("I134").Value = 2228.0715
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6684.2145
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -4149.2145
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3044.6667
$ws.Range("I31").Value = 2225.2222
$ws.Range("J31").Value = 7961.3335
$ws.Range("K31").Value = 2225.2222
$ws.Range("L31").Value = 7961.3335
$ws.Range("M31").Value = -1930.2222
$ws.Range("N31").Value = -8551.333500000001
$ws.Range("H32").Value = 5670
$ws.Range("I32").Value = 5879.8335
$ws.Range("J32").Value = 4411
$ws.Range("K32").Value = 5879.8335
$ws.Range("L32").Value = 4411
$ws.Range("M32").Value = -5563.8335
$ws.Range("N32").Value = -5043
$ws.Range("H34").Value = 3044.6667
$ws.Range("I34").Value = 2225.2222
$ws.Range("J34").Value = 7961.3335
$ws.Range("K34").Value = 2225.2222
$ws.Range("L34").Value = 7961.3335
$ws.Range("M34").Value = -2023.2222
$ws.Range("N34").Value = -8365.333500000001
$ws.Range("H41").Value = 14137.275
$ws.Range("J41").Value = 14356.464
$ws.Range("L41").Value = 14356.464
$ws.Range("N41").Value = -15212.464
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H58").Value = 7326.8237
$ws.Range("I58").Value = 7837.067
$ws.Range("K58").Value = 7837.067
$ws.Range("M58").Value = -7634.067
$ws.Range("H94").Value = 1928.9231
$ws.Range("I94").Value = 1922.625
$ws.Range("K94").Value = 1922.625
$ws.Range("M94").Value = -1471.625
$ws.Range("H105").Value = 3309.9714
$ws.Range("I105").Value = 3389.4736
$ws.Range("J105").Value = 3215.5625
$ws.Range("K105").Value = 3389.4736
$ws.Range("L105").Value = 3215.5625
$ws.Range("M105").Value = -1642.4736
$ws.Range("N105").Value = -6709.5625
$ws.Range("H134").Value = 3246.6365
$ws.Range("I134").Value = 3079.4443
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 9238.332900000001
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -6703.332900000001
$ws.Range("N134").Value = -17067
$ws.Range("H136").Value = 7326.8237
$ws.Range("I136").Value = 7837.067
$ws.Range("K136").Value = 23511.201
$ws.Range("M136").Value = -20961.201

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 107.44444
$ws.Range("J2").Value = 316.8
$ws.Range("L2").Value = 1900.8
$ws.Range("N2").Value = -2126.8
$ws.Range("H12").Value = 455.5
$ws.Range("J12").Value = 171.71428
$ws.Range("L12").Value = 515.14284
$ws.Range("N12").Value = -861.14284
$ws.Range("H17").Value = 117
$ws.Range("I17").Value = 109.416664
$ws.Range("J17").Value = 147.33333
$ws.Range("K17").Value = 328.249992
$ws.Range("L17").Value = 441.99999
$ws.Range("M17").Value = -159.249992
$ws.Range("N17").Value = -779.99999
$ws.Range("H116").Value = 3378.182
$ws.Range("J116").Value = 9666
$ws.Range("L116").Value = 28998
$ws.Range("N116").Value = -35882
$ws.Range("H118").Value = 4058.6667
$ws.Range("I118").Value = 210.5
$ws.Range("J118").Value = 4650.6924
$ws.Range("K118").Value = 631.5
$ws.Range("L118").Value = 13952.0772
$ws.Range("M118").Value = 611.5
$ws.Range("N118").Value = -16438.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3964.6667
$ws.Range("I70").Value = 3964.6667
$ws.Range("K70").Value = 3964.6667
$ws.Range("M70").Value = -3694.6667
$ws.Range("H73").Value = 3964.6667
$ws.Range("I73").Value = 3964.6667
$ws.Range("K73").Value = 3964.6667
$ws.Range("M73").Value = -3028.6667
$ws.Range("H102").Value = 2681.375
$ws.Range("I102").Value = 2783.5
$ws.Range("K102").Value = 2783.5
$ws.Range("M102").Value = -1161.5
$ws.Range("H107").Value = 6179.9
$ws.Range("I107").Value = 6915.1665
$ws.Range("J107").Value = 5077
$ws.Range("K107").Value = 6915.1665
$ws.Range("L107").Value = 5077
$ws.Range("M107").Value = -4995.1665
$ws.Range("N107").Value = -8917
$ws.Range("H122").Value = 1416
$ws.Range("I122").Value = 1416
$ws.Range("K122").Value = 4248
$ws.Range("M122").Value = -1798
$ws.Range("H126").Value = 3499
$ws.Range("I126").Value = 3499
$ws.Range("K126").Value = 10497
$ws.Range("M126").Value = -8027

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1691.8334
$ws.Range("I22").Value = 1958.6666
$ws.Range("J22").Value = 1425
$ws.Range("K22").Value = 1958.6666
$ws.Range("L22").Value = 1425
$ws.Range("M22").Value = -1663.6666
$ws.Range("N22").Value = -2015
$ws.Range("H27").Value = 1691.8334
$ws.Range("I27").Value = 1958.6666
$ws.Range("J27").Value = 1425
$ws.Range("K27").Value = 1958.6666
$ws.Range("L27").Value = 1425
$ws.Range("M27").Value = -1851.6666
$ws.Range("N27").Value = -1639
$ws.Range("H132").Value = 7192.476
$ws.Range("I132").Value = 7474.8335
$ws.Range("K132").Value = 22424.5005
$ws.Range("M132").Value = -19894.5005
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 102780.875
$ws.Range("I62").Value = 4944.4443
$ws.Range("J62").Value = 228570.58
$ws.Range("K62").Value = 4944.4443
$ws.Range("L62").Value = 228570.58
$ws.Range("M62").Value = -4320.4443
$ws.Range("N62").Value = -229818.58
$ws.Range("H65").Value = 102780.875
$ws.Range("I65").Value = 4944.4443
$ws.Range("J65").Value = 228570.58
$ws.Range("K65").Value = 24722.2215
$ws.Range("L65").Value = 1142852.9
$ws.Range("M65").Value = -21602.2215
$ws.Range("N65").Value = -1149092.9
$ws.Range("H81").Value = 2224.3333
$ws.Range("J81").Value = 3015.3333
$ws.Range("L81").Value = 6030.6666
$ws.Range("N81").Value = -8152.6666
$ws.Range("H84").Value = 2224.3333
$ws.Range("J84").Value = 3015.3333
$ws.Range("L84").Value = 30153.333
$ws.Range("N84").Value = -40761.333
$ws.Range("H132").Value = 8905.5
$ws.Range("I132").Value = 9334.929
$ws.Range("J132").Value = 5899.5
$ws.Range("K132").Value = 28004.787
$ws.Range("L132").Value = 17698.5
$ws.Range("M132").Value = -25474.787
$ws.Range("N132").Value = -22758.5
$ws.Range("H136").Value = 2294.037
$ws.Range("I136").Value = 1906.3182
$ws.Range("K136").Value = 5718.9546
$ws.Range("M136").Value = -3168.9546
